# Restore C10 ("Integer min" for rule R30) from 18 to 1, as in the
# source revision being restored.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C10").Value = 1
